# Optimized Tarrival and Tbolus boundary conditions
# Update the optimizer output grid (A1:C16) with the newly converged values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2504.364289323195
$ws.Range("B1").Value = 1661.6143602618283
$ws.Range("C1").Value = 1646.5366247501909
$ws.Range("A2").Value = 2262.5421562324705
$ws.Range("B2").Value = 1496.453505828991
$ws.Range("C2").Value = 1377.738777858284
$ws.Range("A3").Value = 2560.3680356515761
$ws.Range("B3").Value = 1740.6496284113091
$ws.Range("C3").Value = 1579.1624269513977
$ws.Range("A4").Value = 2494.2867125411713
$ws.Range("B4").Value = 1874.5164460460646
$ws.Range("C4").Value = 1886.2354526950858
$ws.Range("A5").Value = 2528.8055215260533
$ws.Range("B5").Value = 1750.1804431919084
$ws.Range("C5").Value = 1760.0114777685721
$ws.Range("A6").Value = 2478.757367532743
$ws.Range("B6").Value = 1811.2349942280243
$ws.Range("C6").Value = 1879.3514180160771
$ws.Range("A7").Value = 2397.2088814474905
$ws.Range("B7").Value = 1847.2451428867341
$ws.Range("C7").Value = 1661.3136774658951
$ws.Range("A8").Value = 2467.9969983093979
$ws.Range("B8").Value = 1930.0311863717509
$ws.Range("C8").Value = 1785.2194460321598
$ws.Range("A9").Value = 2643.8794847871554
$ws.Range("B9").Value = 1943.556967011476
$ws.Range("C9").Value = 1635.6171814922895
$ws.Range("A10").Value = 2394.9100585197939
$ws.Range("B10").Value = 1487.7934969043979
$ws.Range("C10").Value = 1468.8572344874269
$ws.Range("A11").Value = 2163.0420362341006
$ws.Range("B11").Value = 1575.427128416758
$ws.Range("C11").Value = 1407.2591421929676
$ws.Range("A12").Value = 2770.6264727355519
$ws.Range("B12").Value = 2193.9456227505539
$ws.Range("C12").Value = 1875.2102266616878
$ws.Range("A13").Value = 2528.1335108469902
$ws.Range("B13").Value = 1931.1465399730125
$ws.Range("C13").Value = 1747.3254836145011
$ws.Range("A14").Value = 2609.7077477574749
$ws.Range("B14").Value = 2013.3363933907153
$ws.Range("C14").Value = 1762.6916811467902
$ws.Range("A15").Value = 2498.3254311829783
$ws.Range("B15").Value = 1992.2539847176738
$ws.Range("C15").Value = 1845.4903335868835
$ws.Range("A16").Value = 2592.375601129876
$ws.Range("B16").Value = 1796.2307987563786
$ws.Range("C16").Value = 1555.1632149922668
